$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51 with
# refreshed crypto-tracker figures. Values are prefixed with a
# leading apostrophe so Excel stores them as literal text (matching
# the source data, which includes multi-dot thousand separators like
# "26.046.06" alongside plain decimals like "1.008") instead of
# auto-converting plain-looking decimals into floating point numbers.

$ws.Range('D2').Value = "'26.046.06"
$ws.Range('E2').Value = "'  -7.12%  "
$ws.Range('D3').Value = "'1.670.81"
$ws.Range('E3').Value = "'  -4.27%  "
$ws.Range('D4').Value = "'1.008"
$ws.Range('E4').Value = "'  +0.60%  "
$ws.Range('D5').Value = "'217.76"
$ws.Range('E5').Value = "'  -3.63%  "
$ws.Range('D6').Value = "'0.5047"
$ws.Range('E6').Value = "'  -13.03%  "
$ws.Range('D7').Value = "'1.008"
$ws.Range('E7').Value = "'  +0.55%  "
$ws.Range('D8').Value = "'0.2625"
$ws.Range('E8').Value = "'  -3.11%  "
$ws.Range('D9').Value = "'0.06311"
$ws.Range('E9').Value = "'  -4.06%  "
$ws.Range('D10').Value = "'21.31"
$ws.Range('E10').Value = "'  -7.86%  "
$ws.Range('D11').Value = "'0.07367"
$ws.Range('D12').Value = "'1.667.71"
$ws.Range('E12').Value = "'  -4.16%  "
$ws.Range('D13').Value = "'4.527"
$ws.Range('E13').Value = "'  -3.91%  "
$ws.Range('D14').Value = "'0.5730"
$ws.Range('E14').Value = "'  -5.15%  "
$ws.Range('D15').Value = "'1.898.19"
$ws.Range('E15').Value = "'  -4.23%  "
$ws.Range('D16').Value = "'0.000008415"
$ws.Range('E16').Value = "'  -2.68%  "
$ws.Range('D17').Value = "'64.54"
$ws.Range('E17').Value = "'  -12.89%  "
$ws.Range('D18').Value = "'26.154.95"
$ws.Range('E18').Value = "'  -6.66%  "
$ws.Range('D19').Value = "'4.939"
$ws.Range('E19').Value = "'  -7.12%  "
$ws.Range('D20').Value = "'1.007"
$ws.Range('E20').Value = "'  +0.41%  "
$ws.Range('D21').Value = "'10.77"
$ws.Range('E21').Value = "'  -4.28%  "
$ws.Range('D22').Value = "'186.62"
$ws.Range('E22').Value = "'  -8.82%  "
$ws.Range('D23').Value = "'6.163"
$ws.Range('D24').Value = "'1.009"
$ws.Range('E24').Value = "'  +0.54%  "
$ws.Range('D25').Value = "'142.97"
$ws.Range('E25').Value = "'  -5.00%  "
$ws.Range('D26').Value = "'7.599"
$ws.Range('E26').Value = "'  -5.21%  "
$ws.Range('D27').Value = "'0.1166"
$ws.Range('E27').Value = "'  -5.56%  "
$ws.Range('D28').Value = "'15.66"
$ws.Range('E28').Value = "'  -2.69%  "
$ws.Range('D29').Value = "'1.303"
$ws.Range('E29').Value = "'  -6.21%  "
$ws.Range('E30').Value = "'  -5.79%  "
$ws.Range('D31').Value = "'1.323"
$ws.Range('E31').Value = "'  -4.64%  "
$ws.Range('D32').Value = "'3.492"
$ws.Range('E32').Value = "'  -6.66%  "
$ws.Range('D33').Value = "'3.483"
$ws.Range('E33').Value = "'  -6.19%  "
$ws.Range('D34').Value = "'1.662"
$ws.Range('E34').Value = "'  -0.85%  "
$ws.Range('D35').Value = "'1.003"
$ws.Range('E35').Value = "'  -3.15%  "
$ws.Range('D36').Value = "'0.5961"
$ws.Range('E36').Value = "'  -6.10%  "
$ws.Range('D37').Value = "'2.374"
$ws.Range('E37').Value = "'  -3.41%  "
$ws.Range('D38').Value = "'2.646"
$ws.Range('E38').Value = "'  -2.33%  "
$ws.Range('D39').Value = "'0.01597"
$ws.Range('E39').Value = "'  -4.59%  "
$ws.Range('D40').Value = "'1.079.90"
$ws.Range('E40').Value = "'  -4.02%  "
$ws.Range('D41').Value = "'5.908"
$ws.Range('E41').Value = "'  -5.95%  "
$ws.Range('D42').Value = "'0.8572"
$ws.Range('E42').Value = "'  -0.85%  "
$ws.Range('D43').Value = "'1.006"
$ws.Range('E43').Value = "'  +0.04%  "
$ws.Range('D44').Value = "'99.41"
$ws.Range('E44').Value = "'  -0.04%  "
$ws.Range('D45').Value = "'1.821.41"
$ws.Range('D46').Value = "'0.00000000110"
$ws.Range('E46').Value = "'  +3.28%  "
$ws.Range('D47').Value = "'55.77"
$ws.Range('E47').Value = "'  -5.95%  "
$ws.Range('D48').Value = "'1.005"
$ws.Range('E48').Value = "'  +0.66%  "
$ws.Range('D49').Value = "'8.040"
$ws.Range('E49').Value = "'  -2.40%  "
$ws.Range('D50').Value = "'0.4315"
$ws.Range('E50').Value = "'  -2.42%  "
$ws.Range('D51').Value = "'0.05188"
$ws.Range('E51').Value = "'  -3.58%  "

# Drop the "quote prefix" text-override formatting that typing a
# leading apostrophe applies, so the cells' style stays the same
# (unstyled) as every other data cell in the sheet.
$ws.Range('D2:E51').ClearFormats()
